$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Lanita.Kamradt-bda2a6fb"
$ws.Range("B2").Value = "Lanita Kamradt"
$ws.Range("C2").Value = "Leisa.Tatom-3971dd4d"
$ws.Range("D2").Value = 84

$ws.Range("A3").Value = "Bettyann.Dimitt-5808fb09"
$ws.Range("B3").Value = "Bettyann Dimitt"
$ws.Range("C3").Value = "Leisa.Tatom-3971dd4d"
$ws.Range("D3").Value = 81

$ws.Range("A4").Value = "Oralia.Gaekle-77943eed"
$ws.Range("B4").Value = "Oralia Gaekle"
$ws.Range("C4").Value = "Leisa.Tatom-3971dd4d"
$ws.Range("D4").Value = 2

$ws.Range("A5").Value = "Ruben.Busman-0abbd40f"
$ws.Range("B5").Value = "Ruben Busman"
$ws.Range("C5").Value = "Leisa.Tatom-3971dd4d"
$ws.Range("D5").Value = 46

$ws.Range("A6").Value = "Bernetta.Shaske-120785f5"
$ws.Range("B6").Value = "Bernetta Shaske"
$ws.Range("C6").Value = "Leisa.Tatom-3971dd4d"
$ws.Range("D6").Value = 8

$ws.Range("A7").Value = "Leisa.Tatom-3971dd4d"
$ws.Range("B7").Value = "Leisa Tatom"
$ws.Range("D7").Value = 42

$ws.Range("A8").Value = "Renee.Limerick-400a4c02"
$ws.Range("B8").Value = "Renee Limerick"
$ws.Range("C8").Value = "Lanita.Kamradt-bda2a6fb"
$ws.Range("D8").Value = 53

$ws.Range("A9").Value = "Deon.Simcoe-54054eb7"
$ws.Range("B9").Value = "Deon Simcoe"
$ws.Range("C9").Value = "Lanita.Kamradt-bda2a6fb"
$ws.Range("D9").Value = 25

$ws.Range("A10").Value = "Samuel.Zatarain-1f4e7caa"
$ws.Range("B10").Value = "Samuel Zatarain"
$ws.Range("C10").Value = "Lanita.Kamradt-bda2a6fb"
$ws.Range("D10").Value = 86

$ws.Range("A11").Value = "Rosaura.Pajtas-258b3c40"
$ws.Range("B11").Value = "Rosaura Pajtas"
$ws.Range("C11").Value = "Lanita.Kamradt-bda2a6fb"
$ws.Range("D11").Value = 7

$ws.Range("A12").Value = "Phyllis.Macabeo-09608125"
$ws.Range("B12").Value = "Phyllis Macabeo"
$ws.Range("C12").Value = "Lanita.Kamradt-bda2a6fb"
$ws.Range("D12").Value = 34

$ws.Range("A13").Value = "Ruthann.Cruthird-48a47c00"
$ws.Range("B13").Value = "Ruthann Cruthird"
$ws.Range("C13").Value = "Renee.Limerick-400a4c02"
$ws.Range("D13").Value = 87

$ws.Range("A14").Value = "Venus.Viau-58f152ff"
$ws.Range("B14").Value = "Venus Viau"
$ws.Range("C14").Value = "Renee.Limerick-400a4c02"
$ws.Range("D14").Value = 29

$ws.Range("A15").Value = "Rosia.Dobler-ea457b1c"
$ws.Range("B15").Value = "Rosia Dobler"
$ws.Range("C15").Value = "Renee.Limerick-400a4c02"
$ws.Range("D15").Value = 39

$ws.Range("A16").Value = "Mila.Ballinger-53aa6286"
$ws.Range("B16").Value = "Mila Ballinger"
$ws.Range("C16").Value = "Renee.Limerick-400a4c02"
$ws.Range("D16").Value = 83

$ws.Range("A17").Value = "Marylou.Merrit-9744b125"
$ws.Range("B17").Value = "Marylou Merrit"
$ws.Range("C17").Value = "Renee.Limerick-400a4c02"
$ws.Range("D17").Value = 0

